# Generate Report for Handback
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de handback files have been generated: status text, handback
# datetimes, "Latest Target File"/"Latest Handback File" columns + their
# hyperlinks, and a couple of column-width tweaks that came along with the
# report regeneration.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1. "Status" text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shared by the Overview sheet's E/F columns and the per-locale sheets'
#    Status column).
# ---------------------------------------------------------------------
$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

if ($wsOverview.Range("E2").Value2 -eq $oldStatus) { $wsOverview.Range("E2").Value = $newStatus }
if ($wsOverview.Range("F2").Value2 -eq $oldStatus) { $wsOverview.Range("F2").Value = $newStatus }
if ($wsOverview.Range("E3").Value2 -eq $oldStatus) { $wsOverview.Range("E3").Value = $newStatus }
if ($wsOverview.Range("F3").Value2 -eq $oldStatus) { $wsOverview.Range("F3").Value = $newStatus }

if ($wsZhCn.Range("C2").Value2 -eq $oldStatus) { $wsZhCn.Range("C2").Value = $newStatus }
if ($wsZhCn.Range("C3").Value2 -eq $oldStatus) { $wsZhCn.Range("C3").Value = $newStatus }

if ($wsDeDe.Range("C2").Value2 -eq $oldStatus) { $wsDeDe.Range("C2").Value = $newStatus }
if ($wsDeDe.Range("C3").Value2 -eq $oldStatus) { $wsDeDe.Range("C3").Value = $newStatus }

# ---------------------------------------------------------------------
# 2. zh-cn handback finished at 2016-10-26 08:31:32 (was the zero-date
#    placeholder "0001-01-01 00:00:00").
# ---------------------------------------------------------------------
$wsZhCn.Range("K2").Value = "2016-10-26 08:31:32"
$wsZhCn.Range("K3").Value = "2016-10-26 08:31:32"

# ---------------------------------------------------------------------
# 3. de-de handback finished a little later, at 2016-10-26 08:31:48.
# ---------------------------------------------------------------------
$wsDeDe.Range("K2").Value = "2016-10-26 08:31:48"
$wsDeDe.Range("K3").Value = "2016-10-26 08:31:48"

# ---------------------------------------------------------------------
# 4. Fill in "Latest Target File" (I) / "Latest Handback File" (J) for
#    both rows on both locale sheets, with the matching hyperlinks on the
#    "Latest Target File" cells (I2/I3).
# ---------------------------------------------------------------------
$wsZhCn.Range("I2").Value = "45185b63-041e-42c7-80b6-98b651df9ee5.md"
$wsZhCn.Range("J2").Value = "45185b63-041e-42c7-80b6-98b651df9ee5.d6bac2a5174705987fd2065ec770592c303fe6e4.zh-cn.xlf"
$wsZhCn.Range("I3").Value = "b3830289-c780-410f-9b55-a9e2659232bc.md"
$wsZhCn.Range("J3").Value = "b3830289-c780-410f-9b55-a9e2659232bc.e41f4f207819a0d60ca1af9de924095d8a97834f.zh-cn.xlf"

$wsDeDe.Range("I2").Value = "45185b63-041e-42c7-80b6-98b651df9ee5.md"
$wsDeDe.Range("J2").Value = "45185b63-041e-42c7-80b6-98b651df9ee5.d6bac2a5174705987fd2065ec770592c303fe6e4.de-de.xlf"
$wsDeDe.Range("I3").Value = "b3830289-c780-410f-9b55-a9e2659232bc.md"
$wsDeDe.Range("J3").Value = "b3830289-c780-410f-9b55-a9e2659232bc.e41f4f207819a0d60ca1af9de924095d8a97834f.de-de.xlf"

# ---------------------------------------------------------------------
# 5. Re-create the hyperlinks on each locale sheet so the new "Latest
#    Target File" cells (I2/I3) get a link to the same source markdown
#    file as A2/A3, in the same order the original report generator
#    would emit them (A2, I2, A3, I3) so relationship ids line up.
# ---------------------------------------------------------------------
function Set-HandbackHyperlinks($ws) {
    $urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cfe77ee5a4dc650768a299f54012b62f2f25504b/e2e/45185b63-041e-42c7-80b6-98b651df9ee5.md"
    $urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cfe77ee5a4dc650768a299f54012b62f2f25504b/e2e/b3830289-c780-410f-9b55-a9e2659232bc.md"
    $dispA = "45185b63-041e-42c7-80b6-98b651df9ee5.md"
    $dispB = "b3830289-c780-410f-9b55-a9e2659232bc.md"

    $ws.Range("A2").Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $urlA, "", "", $dispA)
    $ws.Hyperlinks.Add($ws.Range("I2"), $urlA, "", "", $dispA)
    $ws.Hyperlinks.Add($ws.Range("A3"), $urlB, "", "", $dispB)
    $ws.Hyperlinks.Add($ws.Range("I3"), $urlB, "", "", $dispB)
}

Set-HandbackHyperlinks $wsZhCn
Set-HandbackHyperlinks $wsDeDe

# ---------------------------------------------------------------------
# 6. Column width tweaks that shipped with the regenerated report: the
#    "datetime" style columns widen from ~17.2 to ~30, and the new
#    target/handback file columns widen to 40 to fit full filenames.
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.17   # E: zh-cn
$wsOverview.Columns.Item(6).ColumnWidth = 29.17   # F: de-de

$wsZhCn.Columns.Item(3).ColumnWidth  = 29.17   # C: Status
$wsZhCn.Columns.Item(9).ColumnWidth  = 39.17   # I: Latest Target File
$wsZhCn.Columns.Item(10).ColumnWidth = 39.17   # J: Latest Handback File

$wsDeDe.Columns.Item(3).ColumnWidth  = 29.17   # C: Status
$wsDeDe.Columns.Item(9).ColumnWidth  = 39.17   # I: Latest Target File
$wsDeDe.Columns.Item(10).ColumnWidth = 39.17   # J: Latest Handback File
